# Update power-law model results for top 100 concepts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: MLE-KS method
$ws.Range("B2").Value = "2.337"
$ws.Range("D2").Value = "10"
$ws.Range("F2").Value = "100"
$ws.Range("H2").Value = "0.055"

# Row 3: Bootstrapping method
$ws.Range("B3").Value = "2.38"
$ws.Range("C3").Value = "0.192"
$ws.Range("D3").Value = "11.602"
$ws.Range("E3").Value = "2.873"
$ws.Range("F3").Value = "87.513"
$ws.Range("G3").Value = "18.628"
$ws.Range("H3").Value = "0.054"
$ws.Range("I3").Value = "0.014"
$ws.Range("J3").Value = "0.438"
